$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated bus-voltage results for the 380 kV case (Case_3_10).
# Each row writes columns B:F and I:N in one shot via a COM SAFEARRAY,
# leaving A (bus index), G (=1) and H (blank) untouched.

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.042599891869828
$rowBF[0,2] = 1.047138268592801
$rowBF[0,3] = 1.050186920802928
$rowBF[0,4] = 1.059784109983779
$ws.Range("B2:F2").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.036005136391901
$rowIN[0,1] = 1.047674844009819
$rowIN[0,2] = 1.049901744698063
$rowIN[0,3] = 1.052941896170506
$rowIN[0,4] = 1.062512670786842
$rowIN[0,5] = 1.019712516440287
$ws.Range("I2:N2").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.043744790363661
$rowBF[0,2] = 1.047997012879938
$rowBF[0,3] = 1.051189517614417
$rowBF[0,4] = 1.060835942231123
$ws.Range("B3:F3").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.036193906685341
$rowIN[0,1] = 1.048465276269072
$rowIN[0,2] = 1.050572150880273
$rowIN[0,3] = 1.053756398689308
$rowIN[0,4] = 1.063378199797882
$rowIN[0,5] = 1.019980809508954
$ws.Range("I3:N3").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.044485795587796
$rowBF[0,2] = 1.048552628430742
$rowBF[0,3] = 1.051838743781196
$rowBF[0,4] = 1.061517028837501
$ws.Range("B4:F4").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.036314650616249
$rowIN[0,1] = 1.048976375819216
$rowIN[0,2] = 1.051005257862792
$rowIN[0,3] = 1.054283302937291
$rowIN[0,4] = 1.063938131321606
$rowIN[0,5] = 1.020154141727477
$ws.Range("I4:N4").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.04479735776766
$rowBF[0,2] = 1.048786197016107
$rowBF[0,3] = 1.052111793000486
$rowBF[0,4] = 1.061803472587022
$ws.Range("B5:F5").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.036365075283553
$rowIN[0,1] = 1.049191155204157
$rowIN[0,2] = 1.051187170521991
$rowIN[0,3] = 1.054504781500243
$rowIN[0,4] = 1.064173496802285
$rowIN[0,5] = 1.020226945546628
$ws.Range("I5:N5").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.04484967300333
$rowBF[0,2] = 1.048825413483925
$rowBF[0,3] = 1.052157645878683
$rowBF[0,4] = 1.061851574488327
$ws.Range("B6:F6").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.036373522093797
$rowIN[0,1] = 1.049227212515044
$rowIN[0,2] = 1.051217704746918
$rowIN[0,3] = 1.054541966848315
$rowIN[0,4] = 1.064213013954716
$rowIN[0,5] = 1.020239165811389
$ws.Range("I6:N6").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.044489958524818
$rowBF[0,2] = 1.048555749433975
$rowBF[0,3] = 1.051842391828335
$rowBF[0,4] = 1.061520855861474
$ws.Range("B7:F7").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.036315325713921
$rowIN[0,1] = 1.048979246053156
$rowIN[0,2] = 1.051007689239236
$rowIN[0,3] = 1.054286262470437
$rowIN[0,4] = 1.063941276403678
$rowIN[0,5] = 1.020155114791365
$ws.Range("I7:N7").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.042986778932081
$rowBF[0,2] = 1.047428495263563
$rowBF[0,3] = 1.050525653832803
$rowBF[0,4] = 1.060139481962121
$ws.Range("B8:F8").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.036069222601633
$rowIN[0,1] = 1.047942049468256
$rowIN[0,2] = 1.050128454504667
$rowIN[0,3] = 1.053217188684748
$rowIN[0,4] = 1.062805205883673
$rowIN[0,5] = 1.019803243469164
$ws.Range("I8:N8").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.040339328641417
$rowBF[0,2] = 1.045441766257062
$rowBF[0,3] = 1.048209067050528
$rowBF[0,4] = 1.057709013906065
$ws.Range("B9:F9").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.035624819538091
$rowIN[0,1] = 1.046111594738264
$rowIN[0,2] = 1.04857384830839
$rowIN[0,3] = 1.051332324982463
$rowIN[0,4] = 1.06080236322046
$rowIN[0,5] = 1.019181126776252
$ws.Range("I9:N9").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.038575221847292
$rowBF[0,2] = 1.044117054428859
$rowBF[0,3] = 1.046667150244259
$rowBF[0,4] = 1.056091188560686
$ws.Range("B10:F10").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.035321340967537
$rowIN[0,1] = 1.044889416029912
$rowIN[0,2] = 1.047533900219311
$rowIN[0,3] = 1.050075064332524
$rowIN[0,4] = 1.059466506745667
$rowIN[0,5] = 1.018764991573799
$ws.Range("I10:N10").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.037811537547348
$rowBF[0,2] = 1.043543389088515
$rowBF[0,3] = 1.046000070120174
$rowBF[0,4] = 1.055391242491938
$ws.Range("B11:F11").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.035188223035246
$rowIN[0,1] = 1.044359752315249
$rowIN[0,2] = 1.047082751479004
$rowIN[0,3] = 1.049530493378124
$rowIN[0,4] = 1.058887917327268
$rowIN[0,5] = 1.018584471119016
$ws.Range("I11:N11").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.037527898060819
$rowBF[0,2] = 1.043330295920303
$rowBF[0,3] = 1.045752373895025
$rowBF[0,4] = 1.055131338860844
$ws.Range("B12:F12").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.035138520327553
$rowIN[0,1] = 1.044162943409447
$rowIN[0,2] = 1.046915047865657
$rowIN[0,3] = 1.049328189952811
$rowIN[0,4] = 1.058672980111381
$rowIN[0,5] = 1.018517368036109
$ws.Range("I12:N12").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.037588738494089
$rowBF[0,2] = 1.043376005502254
$rowBF[0,3] = 1.045805501639383
$rowBF[0,4] = 1.05518708511357
$ws.Range("B13:F13").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.035149193345585
$rowIN[0,1] = 1.044205162679219
$rowIN[0,2] = 1.046951026602235
$rowIN[0,3] = 1.049371585882043
$rowIN[0,4] = 1.058719085931712
$rowIN[0,5] = 1.018531764126313
$ws.Range("I13:N13").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.037788091276728
$rowBF[0,2] = 1.043525774909576
$rowBF[0,3] = 1.045979593698078
$rowBF[0,4] = 1.055369756992138
$ws.Range("B14:F14").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.035184119831682
$rowIN[0,1] = 1.044343485422909
$rowIN[0,2] = 1.047068891634569
$rowIN[0,3] = 1.049513771435558
$rowIN[0,4] = 1.058870151026716
$rowIN[0,5] = 1.01857892537207
$ws.Range("I14:N14").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.037910922635175
$rowBF[0,2] = 1.043618051686687
$rowBF[0,3] = 1.046086869090908
$rowBF[0,4] = 1.055482318758107
$ws.Range("B15:F15").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.035205605169085
$rowIN[0,1] = 1.044428701582768
$rowIN[0,2] = 1.047141495356819
$rowIN[0,3] = 1.049601373260406
$rowIN[0,4] = 1.058963224121353
$rowIN[0,5] = 1.018607976380122
$ws.Range("I15:N15").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.038625908882367
$rowBF[0,2] = 1.044155125533362
$rowBF[0,3] = 1.046711434366438
$rowBF[0,4] = 1.056137653951776
$ws.Range("B16:F16").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.03533013955184
$rowIN[0,1] = 1.044924558514361
$rowIN[0,2] = 1.047563823683329
$rowIN[0,3] = 1.050111202142748
$rowIN[0,4] = 1.059504902551346
$rowIN[0,5] = 1.01877696514506
$ws.Range("I16:N16").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.039074449833618
$rowBF[0,2] = 1.044492002782166
$rowBF[0,3] = 1.047103362877136
$rowBF[0,4] = 1.056548884427954
$ws.Range("B17:F17").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.035407799001576
$rowIN[0,1] = 1.045235475191127
$rowIN[0,2] = 1.047828513006854
$rowIN[0,3] = 1.050430959022602
$rowIN[0,4] = 1.059844641622067
$rowIN[0,5] = 1.018882878674369
$ws.Range("I17:N17").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.039336094155908
$rowBF[0,2] = 1.044688492132261
$rowBF[0,3] = 1.047332024118115
$rowBF[0,4] = 1.056788804462105
$ws.Range("B18:F18").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.035452931476983
$rowIN[0,1] = 1.045416784009161
$rowIN[0,2] = 1.047982820345029
$rowIN[0,3] = 1.050617451525525
$rowIN[0,4] = 1.060042790749393
$rowIN[0,5] = 1.018944624291043
$ws.Range("I18:N18").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.039425311187987
$rowBF[0,2] = 1.044755488993911
$rowBF[0,3] = 1.047410001189178
$rowBF[0,4] = 1.05687062053796
$ws.Range("B19:F19").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.035468292500836
$rowIN[0,1] = 1.045478598240326
$rowIN[0,2] = 1.048035421342507
$rowIN[0,3] = 1.050681037894109
$rowIN[0,4] = 1.060110351952571
$rowIN[0,5] = 1.018965672530704
$ws.Range("I19:N19").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.039026323782504
$rowBF[0,2] = 1.04445585960402
$rowBF[0,3] = 1.047061306868363
$rowBF[0,4] = 1.056504757471468
$ws.Range("B20:F20").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.035399483940321
$rowIN[0,1] = 1.045202121287772
$rowIN[0,2] = 1.047800122773957
$rowIN[0,3] = 1.050396653820051
$rowIN[0,4] = 1.059808192379335
$rowIN[0,5] = 1.018871518460004
$ws.Range("I20:N20").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.037729386117507
$rowBF[0,2] = 1.043481671780441
$rowBF[0,3] = 1.045928325524208
$rowBF[0,4] = 1.055315962262002
$ws.Range("B21:F21").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.035173841939771
$rowIN[0,1] = 1.044302754707843
$rowIN[0,2] = 1.047034186823103
$rowIN[0,3] = 1.0494719020399
$rowIN[0,4] = 1.058825666766611
$rowIN[0,5] = 1.018565038933912
$ws.Range("I21:N21").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.036914104866614
$rowBF[0,2] = 1.042869113176783
$rowBF[0,3] = 1.045216478687336
$rowBF[0,4] = 1.054569025760466
$ws.Range("B22:F22").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.035030486247008
$rowIN[0,1] = 1.043736892416012
$rowIN[0,2] = 1.046551878365684
$rowIN[0,3] = 1.048890326125618
$rowIN[0,4] = 1.058207778808149
$rowIN[0,5] = 1.018372055265615
$ws.Range("I22:N22").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.037346286536153
$rowBF[0,2] = 1.043193846587369
$rowBF[0,3] = 1.045593794474211
$rowBF[0,4] = 1.0549649429754
$ws.Range("B23:F23").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.035106622602828
$rowIN[0,1] = 1.04403690423883
$rowIN[0,2] = 1.046807628801007
$rowIN[0,3] = 1.049198644583706
$rowIN[0,4] = 1.058535345744374
$rowIN[0,5] = 1.01847438684044
$ws.Range("I23:N23").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.03904806981756
$rowBF[0,2] = 1.044472191167189
$rowBF[0,3] = 1.047080309995347
$rowBF[0,4] = 1.056524696369702
$ws.Range("B24:F24").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.035403241668262
$rowIN[0,1] = 1.045217192615593
$rowIN[0,2] = 1.047812951349747
$rowIN[0,3] = 1.050412154914563
$rowIN[0,4] = 1.059824662269147
$rowIN[0,5] = 1.018876651750565
$ws.Range("I24:N24").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.04102360321372
$rowBF[0,2] = 1.04595542446332
$rowBF[0,3] = 1.048807524448167
$rowBF[0,4] = 1.058336910314212
$ws.Range("B25:F25").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.035740979881702
$rowIN[0,1] = 1.046585141357622
$rowIN[0,2] = 1.048976376467942
$rowIN[0,3] = 1.051819728160389
$rowIN[0,4] = 1.06132025710086
$rowIN[0,5] = 1.019342204271673
$ws.Range("I25:N25").Value = $rowIN

